# Insert a new data row at row 103 (pushing existing rows 103:154 down to
# 104:155), then populate the new row with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows(103).Insert()

$ws.Cells.Item(103, 1).Value2 = 4
$ws.Cells.Item(103, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(103, 3).Value = "Los Lagos"
$ws.Cells.Item(103, 4).Value2 = 45072
$ws.Cells.Item(103, 5).Value2 = 10
$ws.Cells.Item(103, 6).Value2 = 100112022
$ws.Cells.Item(103, 7).Value = "Arveja Verde"
$ws.Cells.Item(103, 8).Value = "Perfection"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value2 = 50
$ws.Cells.Item(103, 11).Value2 = 42000
$ws.Cells.Item(103, 12).Value2 = 42000
$ws.Cells.Item(103, 13).Value2 = 42000
$ws.Cells.Item(103, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(103, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(103, 16).Value2 = 1680
$ws.Cells.Item(103, 17).Value2 = 25
$ws.Cells.Item(103, 18).Value = "Hortaliza"
